$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '291.09'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-3.96%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '30.87'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-3.92%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.936'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.16%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07167'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-8.46%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.802'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-10.57%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.660'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-2.18%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.740'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-2.13%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8944'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-2.86%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1650'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-5.91%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07725'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-2.21%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08103'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-5.62%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03051'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-3.49%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1003'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.15%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001498'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.46%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005832'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.00%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.470'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.16%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.080'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-3.61%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.93%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1273'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-3.42%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.040'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-5.51%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1998'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.35%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04513'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.20%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001210'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-1.29%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004007'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-9.95%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.09%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01605'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-7.92%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04384'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-8.30%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007367'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.87%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1305'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-4.24%'
$ws.Range("B43").Value = 'Dexo'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QkL_pl546+dexo-dexo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.007580'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '--%'
$ws.Range("B44").Value = 'CEJI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002019'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-14.45%'
$ws.Range("B45").Value = 'LocalTraders'
$ws.Range("C45").Value = 'https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009251'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-12.61%'
$ws.Range("B46").Value = 'CoinLion'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sot4vgRyjNXek+coinlion-lion'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00005938'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-6.11%'
$ws.Range("B47").Value = 'Kangarootoken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/zkVNkSGwZ3+kangarootoken-gar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000749'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.10%'
$ws.Range("B48").Value = 'BOLO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.247'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '173.92%'
$ws.Range("B49").Value = 'CoinbaseStockToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002998'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-3.34%'
$ws.Range("B50").Value = 'CryptobidCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002098'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.10%'
$ws.Range("B51").Value = 'SpecialPowerGold'
$ws.Range("C51").Value = 'https://coinranking.com/coin/jPTWzmsWb+specialpowergold-spg'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001999'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.10%'
